# Refactor synthetic array /3:
# Swap the "noir" (black) colour-coding symbols/labels for a "bleu" (blue)
# set in the intervention_type legend used throughout the sheet.
#   ⬛ -> 📘
#   🟥 -> 📕
#   🟧 -> 📙
#   🟩 -> 📗
#   noir -> bleu

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

# Use exact whole-cell matches so we only touch the legend cells themselves
# (symbol column and the "noir" label), not any unrelated text that might
# merely contain a similar substring.
$used.Replace("⬛", "📘", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$used.Replace("🟥", "📕", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$used.Replace("🟧", "📙", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$used.Replace("🟩", "📗", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$used.Replace("noir", "bleu", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
